# Apply the documented edit to rest-api.docx:
#  1. Fix typo in the "GET .../hotel/<id>" example request: Verba -> hotel_id
#  2. Resize the API table columns (col1/col2/col3 widths change
#     from 3496/1020/4753 dxa to 2787/1134/5348 dxa, keeping the
#     overall table width the same)

$d = $word.ActiveDocument

# 1) Fix the typo: "Verba" should read "hotel_id" (matches the
#    placeholder used by the other GET/PUT/DELETE hotel examples).
$d.Content.Find.Execute("Verba", $true, $false, $false, $false, $false, `
    $true, 1, $false, "hotel_id", 2) | Out-Null

# 2) Resize the table's three columns. Values are expressed in points
#    (Word COM uses points for Column.Width); dxa/20 = points.
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 2787 / 20   # 139.35 pt (was 174.8 pt / 3496 dxa)
$t.Columns.Item(2).Width = 1134 / 20   # 56.7 pt  (was 51 pt / 1020 dxa)
$t.Columns.Item(3).Width = 5348 / 20   # 267.4 pt (was 237.65 pt / 4753 dxa)
